$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.4801006475809
$ws.Range("C2").Value = 7.950677712841479
$ws.Range("D2").Value = 8.186079208065745
$ws.Range("E2").Value = 12.59484330226064
$ws.Range("F2").Value = 34.55390627995832
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 25.51917288694314
$ws.Range("J2").Value = 9.809917668969984
$ws.Range("L2").Value = 11.31865242984808
$ws.Range("M2").Value = 16.43038510801663
$ws.Range("N2").Value = 18.95874407860837
$ws.Range("O2").Value = 26.52238024649524
$ws.Range("B3").Value = 16.10954675742148
$ws.Range("C3").Value = 7.674397866030133
$ws.Range("D3").Value = 8.19101924888535
$ws.Range("E3").Value = 12.62125984151262
$ws.Range("F3").Value = 34.62114515603839
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 25.61777460339776
$ws.Range("J3").Value = 9.82419855137138
$ws.Range("L3").Value = 11.32182694776376
$ws.Range("M3").Value = 16.35451260464971
$ws.Range("N3").Value = 19.01046269371206
$ws.Range("O3").Value = 26.58807750408752
$ws.Range("B4").Value = 15.8798783727734
$ws.Range("C4").Value = 7.498199367857321
$ws.Range("D4").Value = 8.19480378975701
$ws.Range("E4").Value = 12.63843509017439
$ws.Range("F4").Value = 34.67004136761499
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 25.68294926579048
$ws.Range("J4").Value = 9.83342367168202
$ws.Range("L4").Value = 11.32500197788349
$ws.Range("M4").Value = 16.30979423943714
$ws.Range("N4").Value = 19.0439770176004
$ws.Range("O4").Value = 26.63396365734911
$ws.Range("B5").Value = 15.78587755075004
$ws.Range("C5").Value = 7.424814990294659
$ws.Range("D5").Value = 8.196535475916749
$ws.Range("E5").Value = 12.64567496706303
$ws.Range("F5").Value = 34.6918778680489
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 25.71067273599407
$ws.Range("J5").Value = 9.837298127751929
$ws.Range("L5").Value = 11.32660501688435
$ws.Range("M5").Value = 16.29205349543285
$ws.Range("N5").Value = 19.05807767521958
$ws.Range("O5").Value = 26.65405526740354
$ws.Range("B6").Value = 15.77024799402679
$ws.Range("C6").Value = 7.412536150652832
$ws.Range("D6").Value = 8.196834478677484
$ws.Range("E6").Value = 12.64689170546511
$ws.Range("F6").Value = 34.69561911855071
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 25.71534649604457
$ws.Range("J6").Value = 9.837948442357739
$ws.Range("L6").Value = 11.32688990068157
$ws.Range("M6").Value = 16.28913715210653
$ws.Range("N6").Value = 19.06044588280579
$ws.Range("O6").Value = 26.65747550108941
$ws.Range("B7").Value = 15.87861212487061
$ws.Range("C7").Value = 7.497215990915398
$ws.Range("D7").Value = 8.194826376097774
$ws.Range("E7").Value = 12.63853175381789
$ws.Range("F7").Value = 34.67032812949175
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 25.68331844104791
$ws.Range("J7").Value = 9.833475457297922
$ws.Range("L7").Value = 11.32502234394193
$ws.Range("M7").Value = 16.30955301214829
$ws.Range("N7").Value = 19.04416538775334
$ws.Range("O7").Value = 26.6342289844665
$ws.Range("B8").Value = 16.35285609306536
$ws.Range("C8").Value = 7.856816696285309
$ws.Range("D8").Value = 8.187626914685131
$ws.Range("E8").Value = 12.60375378617608
$ws.Range("F8").Value = 34.57550884316527
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 25.55220880344145
$ws.Range("J8").Value = 9.814747148099569
$ws.Range("L8").Value = 11.31949315226183
$ws.Range("M8").Value = 16.40384442269502
$ws.Range("N8").Value = 18.97621220899619
$ws.Range("O8").Value = 26.54387983664561
$ws.Range("B9").Value = 17.26034225288772
$ws.Range("C9").Value = 8.507319149653423
$ws.Range("D9").Value = 8.179446971072791
$ws.Range("E9").Value = 12.54310833109812
$ws.Range("F9").Value = 34.45007840577223
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 25.33189331479685
$ws.Range("J9").Value = 9.781628826093421
$ws.Range("L9").Value = 11.31833466758912
$ws.Range("M9").Value = 16.60298794424179
$ws.Range("N9").Value = 18.85686616881368
$ws.Range("O9").Value = 26.41082466532656
$ws.Range("B10").Value = 17.90646053998318
$ws.Range("C10").Value = 8.948899688840683
$ws.Range("D10").Value = 8.177025507402112
$ws.Range("E10").Value = 12.50312027911539
$ws.Range("F10").Value = 34.39494853986352
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 25.19249787620111
$ws.Range("J10").Value = 9.759474983064553
$ws.Range("L10").Value = 11.32332830699949
$ws.Range("M10").Value = 16.75717128931018
$ws.Range("N10").Value = 18.77759833899378
$ws.Range("O10").Value = 26.34009120552752
$ws.Range("B11").Value = 18.1945984240857
$ws.Range("C11").Value = 9.141381776478006
$ws.Range("D11").Value = 8.176695925054945
$ws.Range("E11").Value = 12.4859127872325
$ws.Range("F11").Value = 34.37792632886333
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 25.13397197706396
$ws.Range("J11").Value = 9.749865078093846
$ws.Range("L11").Value = 11.32685538155571
$ws.Range("M11").Value = 16.82884198508048
$ws.Range("N11").Value = 18.74335132351011
$ws.Range("O11").Value = 26.31380250129345
$ws.Range("B12").Value = 18.30277651075565
$ws.Range("C12").Value = 9.213026904177338
$ws.Range("D12").Value = 8.176681471940341
$ws.Range("E12").Value = 12.47953756306037
$ws.Range("F12").Value = 34.37263950004062
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 25.11251322630933
$ws.Range("J12").Value = 9.746293014973666
$ws.Range("L12").Value = 11.3283702227124
$ws.Range("M12").Value = 16.85618625334682
$ws.Range("N12").Value = 18.73064250295626
$ws.Range("O12").Value = 26.30469574864021
$ws.Range("B13").Value = 18.27952148022231
$ws.Range("C13").Value = 9.197652685109849
$ws.Range("D13").Value = 8.176679687390875
$ws.Range("E13").Value = 12.48090432412451
$ws.Range("F13").Value = 34.373726558459
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 25.11710343460724
$ws.Range("J13").Value = 9.747059347566656
$ws.Range("L13").Value = 11.3280360256435
$ws.Range("M13").Value = 16.85028835742246
$ws.Range("N13").Value = 18.73336803400679
$ws.Range("O13").Value = 26.30661930590374
$ws.Range("B14").Value = 18.20351756442657
$ws.Range("C14").Value = 9.147301180600175
$ws.Range("D14").Value = 8.176692528190257
$ws.Range("E14").Value = 12.48538547319092
$ws.Range("F14").Value = 34.37746814671706
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 25.13219244205261
$ws.Range("J14").Value = 9.749569861236136
$ws.Range("L14").Value = 11.32697642450724
$ws.Range("M14").Value = 16.83108761695792
$ws.Range("N14").Value = 18.74230055947236
$ws.Range("O14").Value = 26.31303627643417
$ws.Range("B15").Value = 18.15683851540823
$ws.Range("C15").Value = 9.116296470203226
$ws.Range("D15").Value = 8.176714744893053
$ws.Range("E15").Value = 12.48814863790227
$ws.Range("F15").Value = 34.37991093535053
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 25.14152657700893
$ws.Range("J15").Value = 9.751116339870778
$ws.Range("L15").Value = 11.32635068699895
$ws.Range("M15").Value = 16.81935270847461
$ws.Range("N15").Value = 18.74780579613981
$ws.Range("O15").Value = 26.31707735709649
$ws.Range("B16").Value = 17.88750479056151
$ws.Range("C16").Value = 8.936148371525205
$ws.Range("D16").Value = 8.177062524780695
$ws.Range("E16").Value = 12.50426457109216
$ws.Range("F16").Value = 34.39622316446749
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 25.19642109529023
$ws.Range("J16").Value = 9.760112404650437
$ws.Range("L16").Value = 11.32312295612471
$ws.Range("M16").Value = 16.75251689257629
$ws.Range("N16").Value = 18.77987285923339
$ws.Range("O16").Value = 26.34192785465576
$ws.Range("B17").Value = 17.72072100722218
$ws.Range("C17").Value = 8.823455635402027
$ws.Range("D17").Value = 8.177473135347348
$ws.Range("E17").Value = 12.51440264656658
$ws.Range("F17").Value = 34.40829428570977
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 25.23134929699854
$ws.Range("J17").Value = 9.765750850163951
$ws.Range("L17").Value = 11.32146348186074
$ws.Range("M17").Value = 16.71189624650027
$ws.Range("N17").Value = 18.80000854652756
$ws.Range("O17").Value = 26.35868210609473
$ws.Range("B18").Value = 17.62425252193902
$ws.Range("C18").Value = 8.75784941373068
$ws.Range("D18").Value = 8.177782008301348
$ws.Range("E18").Value = 12.52032638018668
$ws.Range("F18").Value = 34.41599560700465
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 25.25189887726389
$ws.Range("J18").Value = 9.769038001849305
$ws.Range("L18").Value = 11.32062724749898
$ws.Range("M18").Value = 16.68867750408127
$ws.Range("N18").Value = 18.81176069349647
$ws.Range("O18").Value = 26.36887287222874
$ws.Range("B19").Value = 17.59150061982835
$ws.Range("C19").Value = 8.735502068244337
$ws.Range("D19").Value = 8.177899093270044
$ws.Range("E19").Value = 12.52234797118873
$ws.Range("F19").Value = 34.41873335903223
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 25.25893554825377
$ws.Range("J19").Value = 9.770158552564002
$ws.Range("L19").Value = 11.3203644616865
$ws.Range("M19").Value = 16.68084145882094
$ws.Range("N19").Value = 18.81576910033511
$ws.Range("O19").Value = 26.37241841909717
$ws.Range("B20").Value = 17.73853195620808
$ws.Range("C20").Value = 8.835533839476808
$ws.Range("D20").Value = 8.177421905570984
$ws.Range("E20").Value = 12.51331385310953
$ws.Range("F20").Value = 34.40693080369515
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 25.22758353404579
$ws.Range("J20").Value = 9.765146069408125
$ws.Range("L20").Value = 11.32162790647966
$ws.Range("M20").Value = 16.71620547424738
$ws.Range("N20").Value = 18.79784741450084
$ws.Range("O20").Value = 26.35684122193127
$ws.Range("B21").Value = 18.22586786433903
$ws.Range("C21").Value = 9.162124657036049
$ws.Range("D21").Value = 8.17668576693128
$ws.Range("E21").Value = 12.48406543172405
$ws.Range("F21").Value = 34.37633769188184
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 25.12774132621182
$ws.Range("J21").Value = 9.748830646526793
$ws.Range("L21").Value = 11.32728280204658
$ws.Range("M21").Value = 16.83672192998529
$ws.Range("N21").Value = 18.73966981633939
$ws.Range("O21").Value = 26.31112842610538
$ws.Range("B22").Value = 18.53888713759104
$ws.Range("C22").Value = 9.368308050046318
$ws.Range("D22").Value = 8.176847556117599
$ws.Range("E22").Value = 12.46577081682503
$ws.Range("F22").Value = 34.3630993703077
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 25.06659114075891
$ws.Range("J22").Value = 9.738557961955163
$ws.Range("L22").Value = 11.33202254374174
$ws.Range("M22").Value = 16.91666854816273
$ws.Range("N22").Value = 18.703161224362
$ws.Range("O22").Value = 26.28619696778142
$ws.Range("B23").Value = 18.37235620694562
$ws.Range("C23").Value = 9.258939427579241
$ws.Range("D23").Value = 8.176702603290984
$ws.Range("E23").Value = 12.47546005130764
$ws.Range("F23").Value = 34.36954669571043
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 25.0988523803571
$ws.Range("J23").Value = 9.744005063996051
$ws.Range("L23").Value = 11.32939778642835
$ws.Range("M23").Value = 16.87389677606098
$ws.Range("N23").Value = 18.7225082999451
$ws.Range("O23").Value = 26.29905051954461
$ws.Range("B24").Value = 17.73048143197717
$ws.Range("C24").Value = 8.83007582748939
$ws.Range("D24").Value = 8.177444839728247
$ws.Range("E24").Value = 12.51380580000118
$ws.Range("F24").Value = 34.40754486202021
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 25.22928457498141
$ws.Range("J24").Value = 9.765419348930147
$ws.Range("L24").Value = 11.32155320301719
$ws.Range("M24").Value = 16.7142568518581
$ws.Range("N24").Value = 18.79882391438779
$ws.Range("O24").Value = 26.35767174599766
$ws.Range("B25").Value = 17.01799005618277
$ws.Range("C25").Value = 8.337530140963628
$ws.Range("D25").Value = 8.181027323146065
$ws.Range("E25").Value = 12.55870966525313
$ws.Range("F25").Value = 34.47751659659484
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 25.38755157784143
$ws.Range("J25").Value = 9.790204171981324
$ws.Range("L25").Value = 11.31761713263728
$ws.Range("M25").Value = 16.54767158026258
$ws.Range("N25").Value = 18.88766998463011
$ws.Range("O25").Value = 26.44208238668772

$wb.Save()